# "renamed list names | databases"
#
# The source workbook has a data worksheet "Ведомость" (plus a second,
# chart-only sheet "Ведомость(Диаграмма)" that this COM surface does not
# expose as a renameable Sheets/Worksheets member). Rename the data
# worksheet to "Относительные ссылки" and keep every formula/chart
# reference in sync with the new name.

$wb = $excel.ActiveWorkbook

# Locate the data worksheet robustly: prefer the known original name, but
# fall back to "the first/only worksheet" if that lookup ever fails.
$oldName = "Ведомость"
$newName = "Относительные ссылки"

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq $oldName) { $ws = $sheet }
}
if ($ws -eq $null) { $ws = $wb.Worksheets.Item(1) }

$ws.Name = $newName

# The engine does not cascade a worksheet rename into cached chart
# formulas, so fix up every series on every chart embedded on this sheet
# so the cat/val references keep pointing at '<newName>'!... instead of
# the stale '<oldName>'!... range.
for ($ci = 1; $ci -le $ws.ChartObjects().Count; $ci++) {
    $chart = $ws.ChartObjects().Item($ci).Chart
    for ($si = 1; $si -le $chart.SeriesCollection().Count; $si++) {
        $series = $chart.SeriesCollection().Item($si)
        $series.Formula = $series.Formula.Replace("'" + $oldName + "'!", "'" + $newName + "'!").Replace($oldName + "!", "'" + $newName + "'!")
    }
}
